$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: new data row -------------------------------------------------
# A3 "test@email.com" -- same look as A2/B2 (bordered data-row style)
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = "test@email.com"

# B3 "welcome123" -- bordered, but default (un-bolded, non-Arial) font
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B3").Value = "welcome123"
$ws.Range("B3").Font.Name = "Calibri"
$ws.Range("B3").Font.Size = 11
$ws.Range("B3").Font.Bold = $false
$ws.Range("B3").Font.ThemeFont = 1

$excel.CutCopyMode = 0

# Login-page locator: a clickable hyperlink on the e-mail address
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:test@email.com")

# Hyperlinks.Add recolors/underlines the cell -- restore the plain data-row look
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$wb.Styles.Item("Hyperlink").Delete()

# Column B now holds real content -- size it to fit
$ws.Columns.Item(2).AutoFit()

# Mirror the saved selection state recorded in the workbook
$ws.Range("B2").Select() | Out-Null

$wb.Save() | Out-Null
